$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108, shifting existing rows 108-111 down to 109-112
$ws.Rows.Item(108).Insert()

# Populate the new row 108 with the new data
$ws.Cells.Item(108, 1).Value = 2
$ws.Cells.Item(108, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(108, 3).Value = "Coquimbo"
$ws.Cells.Item(108, 4).Value = 44595
$ws.Cells.Item(108, 5).Value = 4
$ws.Cells.Item(108, 6).Value = 100112024
$ws.Cells.Item(108, 7).Value = "Choclo"
$ws.Cells.Item(108, 8).Value = "Choclero"
$ws.Cells.Item(108, 9).Value = "Primera"
$ws.Cells.Item(108, 10).Value = 50000
$ws.Cells.Item(108, 11).Value = 150
$ws.Cells.Item(108, 12).Value = 180
$ws.Cells.Item(108, 13).Value = 165
$ws.Cells.Item(108, 14).Value = "$/unidad"
$ws.Cells.Item(108, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(108, 16).Value = 165
$ws.Cells.Item(108, 17).Value = 1
$ws.Cells.Item(108, 18).Value = "Hortaliza"
